# Insert a new row at row 14 (pushes existing rows 14-83 down to 15-84)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(14).Insert()

# Populate the new row 14 with the new data record
$ws.Cells.Item(14, 1).Value = 1
$ws.Cells.Item(14, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(14, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(14, 4).Value = "2022-04-26"
$ws.Cells.Item(14, 5).Value = 15
$ws.Cells.Item(14, 6).Value = 100112038
$ws.Cells.Item(14, 7).Value = "Cebollín baby"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 300
$ws.Cells.Item(14, 11).Value = 2300
$ws.Cells.Item(14, 12).Value = 2500
$ws.Cells.Item(14, 13).Value = 2400
$ws.Cells.Item(14, 14).Value = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(14, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(14, 16).Value = 1200
$ws.Cells.Item(14, 17).Value = 2
$ws.Cells.Item(14, 18).Value = "Hortaliza"
